# Update the cryptos list (prices / 1h volume %) per the GitHub Actions refresh.
# Note: several "Price" values are plain numeric-looking text (e.g. "597.08")
# that must stay as text (matching the original inlineStr cells) rather than
# being auto-converted to numbers by Excel's input parser, so we force those
# cells to Text format ("@") before assigning them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.681.86'
$ws.Range('E2').Value = '  +0.56%  '
$ws.Range('D3').Value = '3.797.45'
$ws.Range('E3').Value = '  +1.07%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.08'
$ws.Range('E5').Value = '  +0.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.13'
$ws.Range('E6').Value = '  +0.82%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  +0.68%  '
$ws.Range('E9').Value = '  +1.52%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.31'
$ws.Range('E10').Value = '  -0.73%  '
$ws.Range('E11').Value = '  +0.39%  '
$ws.Range('E12').Value = '  -0.18%  '
$ws.Range('E13').Value = '  +0.49%  '
$ws.Range('D14').Value = '4.439.63'
$ws.Range('E14').Value = '  +1.15%  '
$ws.Range('D15').Value = '3.783.06'
$ws.Range('E15').Value = '  +1.36%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '18.59'
$ws.Range('E16').Value = '  +5.07%  '
$ws.Range('D17').Value = '67.646.83'
$ws.Range('E17').Value = '  +0.57%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.09'
$ws.Range('E18').Value = '  +2.37%  '
$ws.Range('E19').Value = '  +0.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '460.68'
$ws.Range('E20').Value = '  +1.22%  '
$ws.Range('E21').Value = '  -5.49%  '
$ws.Range('E22').Value = '  +1.17%  '
$ws.Range('E23').Value = '  +0.39%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.40'
$ws.Range('E24').Value = '  +0.47%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.09'
$ws.Range('E25').Value = '  +2.60%  '
$ws.Range('E26').Value = '  -0.17%  '
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('E28').Value = '  +0.63%  '
$ws.Range('D29').Value = '3.940.33'
$ws.Range('E29').Value = '  +0.93%  '
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('E31').Value = '  +3.05%  '
$ws.Range('E32').Value = '  +1.45%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '29.56'
$ws.Range('E33').Value = '  -0.10%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.07'
$ws.Range('E35').Value = '  -1.19%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.1000'
$ws.Range('E36').Value = '  +0.89%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.37'
$ws.Range('E37').Value = '  +1.91%  '
$ws.Range('E38').Value = '  +0.52%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  +1.07%  '
$ws.Range('E40').Value = '  +0.94%  '
$ws.Range('E41').Value = '  -0.07%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '48.11'
$ws.Range('E43').Value = '  +3.04%  '
$ws.Range('E44').Value = '  +1.38%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '43.05'
$ws.Range('E45').Value = '  -0.92%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '147.67'
$ws.Range('E47').Value = '  +0.33%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '27.19'
$ws.Range('E48').Value = '  +8.77%  '
$ws.Range('B49').Value = 'ONDO'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.36'
$ws.Range('E49').Value = '  +11.10%  '
$ws.Range('B50').Value = 'Bittensor'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '394.47'
$ws.Range('E50').Value = '  +1.69%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.85'
$ws.Range('E51').Value = '  +1.30%  '
